# Update automàtic: dades i banners [2026-02-14 17:50]
# Refreshes DATA_EXTRACCIO timestamps and the latest meteo.cat readings
# for each station row on the Dades_Meteo sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-14 17:48:25"
$ws.Range("I2").Value = "31.0 mm"
$ws.Range("N2").Value = "-2.4 °C 17:29 TU"
$ws.Range("E3").Value = "2026-02-14 17:48:28"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "96%"
$ws.Range("I3").Value = "12.7 mm"
$ws.Range("K3").Value = "5.4 MJ/m2"
$ws.Range("N3").Value = "-6.1 °C 17:29 TU"
$ws.Range("O3").Value = "-4.9 °C"
$ws.Range("E4").Value = "2026-02-14 17:48:30"
$ws.Range("J4").Value = "995.3 hPa"
$ws.Range("E5").Value = "2026-02-14 17:48:32"
$ws.Range("I5").Value = "19.3 mm"
$ws.Range("L5").Value = "50.8 km/h - 329º 17:29 TU"
$ws.Range("N5").Value = "-6.0 °C 17:13 TU"
$ws.Range("E6").Value = "2026-02-14 17:48:35"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "81%"
$ws.Range("J6").Value = "995.4 hPa"
$ws.Range("O6").Value = "10.0 °C"
$ws.Range("E7").Value = "2026-02-14 17:48:38"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "56%"
$ws.Range("J7").Value = "995.6 hPa"
$ws.Range("O7").Value = "13.0 °C"
$ws.Range("E8").Value = "2026-02-14 17:48:40"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "67%"
$ws.Range("J8").Value = "995.3 hPa"
$ws.Range("O8").Value = "9.5 °C"
$ws.Range("E9").Value = "2026-02-14 17:48:43"
$ws.Range("E10").Value = "2026-02-14 17:48:45"
$ws.Range("L10").Value = "31.0 km/h - 39º 17:20 TU"
$ws.Range("O10").Value = "10.0 °C"
$ws.Range("E11").Value = "2026-02-14 17:48:48"
$ws.Range("E12").Value = "2026-02-14 17:48:50"
$ws.Range("E13").Value = "2026-02-14 17:48:52"
$ws.Range("J13").Value = "997.7 hPa"
$ws.Range("E14").Value = "2026-02-14 17:48:55"
$ws.Range("K14").Value = "14.0 MJ/m2"
$ws.Range("E15").Value = "2026-02-14 17:48:57"
$ws.Range("E16").Value = "2026-02-14 17:49:00"
$ws.Range("K16").Value = "8.9 MJ/m2"
$ws.Range("N16").Value = "-8.0 °C 17:29 TU"
$ws.Range("E17").Value = "2026-02-14 17:49:02"
$ws.Range("E18").Value = "2026-02-14 17:49:05"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "79%"
$ws.Range("J18").Value = "995.6 hPa"
$ws.Range("E19").Value = "2026-02-14 17:49:07"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "82%"
$ws.Range("E20").Value = "2026-02-14 17:49:10"
$ws.Range("I20").Value = "0.1 mm"
$ws.Range("N20").Value = "-6.1 °C 17:28 TU"
$ws.Range("E21").Value = "2026-02-14 17:49:12"
$ws.Range("J21").Value = "997.6 hPa"
$ws.Range("E22").Value = "2026-02-14 17:49:14"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "86%"
$ws.Range("N22").Value = "-8.4 °C 17:02 TU"
$ws.Range("E23").Value = "2026-02-14 17:49:17"
$ws.Range("I23").Value = "34.5 mm"
$ws.Range("N23").Value = "-7.2 °C 17:22 TU"
$ws.Range("E24").Value = "2026-02-14 17:49:19"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "68%"
$ws.Range("J24").Value = "999.4 hPa"
$ws.Range("K24").Value = "14.6 MJ/m2"
$ws.Range("E25").Value = "2026-02-14 17:49:22"
$ws.Range("I25").Value = "6.7 mm"
$ws.Range("N25").Value = "-6.0 °C 17:21 TU"
$ws.Range("O25").Value = "-4.4 °C"
$ws.Range("E26").Value = "2026-02-14 17:49:24"
$ws.Range("E27").Value = "2026-02-14 17:49:27"
$ws.Range("N27").Value = "-3.6 °C 17:29 TU"
$ws.Range("E28").Value = "2026-02-14 17:49:29"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "74%"
$ws.Range("J28").Value = "995.3 hPa"
$ws.Range("O28").Value = "8.9 °C"
$ws.Range("E29").Value = "2026-02-14 17:49:31"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "66%"
$ws.Range("K29").Value = "8.7 MJ/m2"
$ws.Range("E30").Value = "2026-02-14 17:49:34"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "58%"
$ws.Range("J30").Value = "995.2 hPa"
$ws.Range("E31").Value = "2026-02-14 17:49:36"
$ws.Range("J31").Value = "994.4 hPa"
$ws.Range("E32").Value = "2026-02-14 17:49:39"
$ws.Range("K32").Value = "11.1 MJ/m2"
$ws.Range("O32").Value = "4.3 °C"
$ws.Range("E33").Value = "2026-02-14 17:49:41"
$ws.Range("J33").Value = "997.0 hPa"
$ws.Range("E34").Value = "2026-02-14 17:49:44"
$ws.Range("I34").Value = "3.0 mm"
$ws.Range("N34").Value = "-3.6 °C 17:27 TU"
$ws.Range("O34").Value = "-2.0 °C"
$ws.Range("E35").Value = "2026-02-14 17:49:46"
$ws.Range("J35").Value = "1002.2 hPa"
$ws.Range("E36").Value = "2026-02-14 17:49:49"
$ws.Range("J36").Value = "996.0 hPa"
$ws.Range("E37").Value = "2026-02-14 17:49:51"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "70%"
$ws.Range("J37").Value = "996.2 hPa"
$ws.Range("L37").Value = "43.6 km/h - 284º 17:08 TU"
$ws.Range("O37").Value = "6.7 °C"
$ws.Range("E38").Value = "2026-02-14 17:49:54"
$ws.Range("L38").Value = "20.2 km/h - 73º 17:09 TU"
$ws.Range("E39").Value = "2026-02-14 17:49:56"
$ws.Range("I39").Value = "11.0 mm"
$ws.Range("N39").Value = "-7.2 °C 17:28 TU"
$ws.Range("E40").Value = "2026-02-14 17:49:59"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "65%"
$ws.Range("J40").Value = "998.0 hPa"
$ws.Range("O40").Value = "7.3 °C"
$ws.Range("E41").Value = "2026-02-14 17:50:01"
$ws.Range("J41").Value = "997.3 hPa"
$ws.Range("O41").Value = "13.3 °C"
$ws.Range("E42").Value = "2026-02-14 17:50:04"
$ws.Range("E43").Value = "2026-02-14 17:50:06"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "71%"
$ws.Range("O43").Value = "9.0 °C"
$ws.Range("E44").Value = "2026-02-14 17:50:09"
$ws.Range("I44").Value = "33.4 mm"
$ws.Range("N44").Value = "-6.3 °C 17:02 TU"
$ws.Range("O44").Value = "-5.2 °C"
$ws.Range("E45").Value = "2026-02-14 17:50:11"
$ws.Range("I45").Value = "13.0 mm"
$ws.Range("J45").Value = "1004.3 hPa"
$ws.Range("K45").Value = "2.1 MJ/m2"
$ws.Range("E46").Value = "2026-02-14 17:50:14"
